# Updated symbol list on Thu Dec 15 19:25:37 UTC 2022 with GitHub Actions
# Applies the updated "Price" column (column D) values to the cryptos sheet.
# Values are stored as text (matching the workbook's inlineStr convention),
# so a leading apostrophe is used to prevent Excel from auto-converting the
# numeric-looking strings into floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "260.64"
    "D4"  = "6.202"
    "D5"  = "0.06088"
    "D6"  = "3.515"
    "D7"  = "6.704"
    "D8"  = "1.357"
    "D9"  = "0.7982"
    "D11" = "0.08076"
    "D12" = "0.03339"
    "D13" = "0.03121"
    "D14" = "0.09277"
    "D15" = "3.942"
    "D16" = "0.001707"
    "D17" = "0.04811"
    "D18" = "0.0006159"
    "D19" = "0.006188"
    "D21" = "0.003392"
    "D22" = "0.0001500"
    "D23" = "3.694"
    "D24" = "2.292"
    "D25" = "0.3358"
    "D26" = "0.1183"
    "D27" = "0.0006164"
    "D40" = "0.04602"
    "D41" = "0.007138"
    "D42" = "0.003899"
    "D43" = "0.1118"
    "D45" = "0.002969"
    "D46" = "0.00006009"
    "D48" = "0.7498"
    "D49" = "0.1277"
    "D50" = "0.00001500"
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = "'" + $updates[$cell]
}
